$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 14, pushing the existing row 14 (and below)
# down by one. This turns the old row 14 into row 15, and the old row 15
# into row 16.
$ws.Rows("14:14").Insert()

# Populate the newly inserted row 14 with the latest weekly record.
$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Vega Modelo de Temuco"
$ws.Range("C14").Value = "La Araucanía"
$ws.Range("D14").Value = 44508
$ws.Range("D14").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 100112017
$ws.Range("G14").Value = "Ramas de apio"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 30
$ws.Range("K14").Value = 4000
$ws.Range("L14").Value = 4000
$ws.Range("M14").Value = 4000
$ws.Range("N14").Value = "$/paquete"
$ws.Range("O14").Value = "Región de La Araucanía"
$ws.Range("P14").Value = 4000
$ws.Range("Q14").Value = 1
$ws.Range("R14").Value = "Hortaliza"
